# The presentation currently applies the "Integral" (Red Violet) theme to
# the slide master (backed by ppt/theme/theme2.xml) while ppt/theme/theme1.xml
# (used by the notes master) carries the default "Office Theme" colors.
#
# The target revision swaps the two: the slides should go back to showing
# the default Office color palette. We reach the theme that backs the
# slide master through the Slide.ThemeColorScheme collection (12 slots,
# same order as a:clrScheme: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# and push the standard Office theme RGB values into it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0        # dk1      #000000
$tcs.Item(2).RGB  = 16777215 # lt1      #FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      #44546A
$tcs.Item(4).RGB  = 15132391 # lt2      #E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  #5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  #ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  #A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  #FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  #4472C4
$tcs.Item(10).RGB = 4697456  # accent6  #70AD47
$tcs.Item(11).RGB = 12673797 # hlink    #0563C1
$tcs.Item(12).RGB = 7491477  # folHlink #954F72
